# "Ändrade funktionen clean data"
# Update the label text in A7 from "Klassificeringsstruktur" to
# "KlassificeringsstrukturText", and move the active selection to B11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "KlassificeringsstrukturText"

$ws.Range("B11").Select()
